# Auto-generated edit script applying cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Thu Aug 15 11:19:54 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.416.40'
$ws.Range('E2').Value = '  -3.80%  '
$ws.Range('D3').Value = '2.617.58'
$ws.Range('E3').Value = '  -3.45%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.18'
$ws.Range('E5').Value = '  -1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.32'
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -1.77%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.58'
$ws.Range('E9').Value = '  -7.13%  '
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.335'
$ws.Range('E11').Value = '  -1.79%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').Value = '3.078.70'
$ws.Range('E13').Value = '  -3.51%  '
$ws.Range('D14').Value = '58.378.16'
$ws.Range('E14').Value = '  -3.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.94'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.652.43'
$ws.Range('E16').Value = '  -2.61%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '336.94'
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.40'
$ws.Range('E19').Value = '  -2.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.37'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.27'
$ws.Range('E21').Value = '  -2.41%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.88'
$ws.Range('E23').Value = '  +2.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.414'
$ws.Range('E24').Value = '  -1.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.166'
$ws.Range('E25').Value = '  -2.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.12'
$ws.Range('E27').Value = '  -2.68%  '
$ws.Range('D28').Value = '0.0₃0789'
$ws.Range('E28').Value = '  -4.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.52'
$ws.Range('E29').Value = '  -3.38%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  -0.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.74'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '150.08'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.09'
$ws.Range('E34').Value = '  -3.87%  '
$ws.Range('E35').Value = '  -3.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.891'
$ws.Range('E36').Value = '  -3.11%  '
$ws.Range('E37').Value = '  -5.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.24'
$ws.Range('E38').Value = '  -2.87%  '
$ws.Range('E39').Value = '  -6.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.62'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.602'
$ws.Range('E42').Value = '  -3.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0970'
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '267.51'
$ws.Range('E44').Value = '  -4.78%  '
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.09'
$ws.Range('E46').Value = '  -5.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0529'
$ws.Range('E47').Value = '  -2.38%  '
$ws.Range('D48').Value = '2.031.03'
$ws.Range('E48').Value = '  -3.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0227'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('E50').Value = '  -7.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.26'
$ws.Range('E51').Value = '  -5.34%  '
